# Release-Notes.xlsx update
# "Microsoft Defender for Cloud - v1" was refreshed (new timestamp) and, since
# the Folder Inventory sheet is kept sorted by Last Updated (desc), it moves
# from the bottom of the list back up to the top. Every other folder's row
# shifts down by one to make room, and the stale duplicate entry that falls
# off the shifted block (the old dated "Microsoft Defender for Cloud - v1"
# row) is removed so the sheet keeps the same 74 data rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Folder Inventory sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Folder Inventory")

# Insert a fresh row right under the header, pushing rows 2..75 to 3..76.
$ws.Rows("2:2").Insert()
# The inserted row inherits the header's bold/border style by default -
# strip that back to the plain formatting the rest of the data rows use.
$ws.Rows("2:2").ClearFormats()

$ws.Range("A2").Value = "Microsoft Defender for Cloud - v1"
$ws.Range("B2").Value = "Microsoft Defender for Cloud - v1"
$ws.Range("C2").Value = "2025-06-16 17:12:47 +0530"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "Root"

# The old ("Microsoft Defender for Cloud - v1", 2025-06-05 19:52:09 +0530)
# row has now been pushed down to row 39 - it's superseded by the refreshed
# entry above, so drop it and let everything below settle back up.
$ws.Rows("39:39").Delete()

# ---------------------------------------------------------------------
# 2) Metadata sheet
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2025-06-16 11:43:19 UTC"
# "Workflow Run" is stored as text ("1"/"2"), not a number - force text via a
# leading apostrophe then strip the resulting quote-prefix formatting so the
# cell keeps the plain (unstyled) look the rest of the sheet uses.
$meta.Range("B5").Value = "'2"
$meta.Range("B5").ClearFormats()

# ---------------------------------------------------------------------
# 3) Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = "2025-06-16 17:12:47 +0530"
